# Updates Price (column D) values on Sheet1 per the "Updated symbol list" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell, expected prior price text, new price text.
# Prices are stored as text (not numbers) in this sheet, so we force the
# cell to Text format before writing the value -- otherwise a numeric-looking
# string like "272.46" would be auto-converted to a floating point number.
$changes = @(
    @{Cell="D2"; Old="272.48"; New="272.46"},
    @{Cell="D4"; Old="6.355"; New="6.352"},
    @{Cell="D5"; Old="0.06335"; New="0.06341"},
    @{Cell="D6"; Old="3.662"; New="3.668"},
    @{Cell="D7"; Old="6.768"; New="6.780"},
    @{Cell="D8"; Old="1.402"; New="1.414"},
    @{Cell="D9"; Old="0.8364"; New="0.8382"},
    @{Cell="D10"; Old="0.1628"; New="0.1630"},
    @{Cell="D11"; Old="0.08392"; New="0.08416"},
    @{Cell="D12"; Old="0.03414"; New="0.03415"},
    @{Cell="D13"; Old="0.03158"; New="0.03160"},
    @{Cell="D14"; Old="0.09307"; New="0.09311"},
    @{Cell="D15"; Old="3.897"; New="3.903"},
    @{Cell="D16"; Old="0.001723"; New="0.001733"},
    @{Cell="D17"; Old="0.04869"; New="0.04874"},
    @{Cell="D18"; Old="0.006204"; New="0.006203"},
    @{Cell="D19"; Old="0.005505"; New="0.005503"},
    @{Cell="D21"; Old="0.0001496"; New="0.0001497"},
    @{Cell="D22"; Old="3.744"; New="3.737"},
    @{Cell="D23"; Old="2.302"; New="2.344"},
    @{Cell="D25"; Old="0.3338"; New="0.3346"},
    @{Cell="D27"; Old="0.0002674"; New="0.0002675"},
    @{Cell="D40"; Old="0.04695"; New="0.04697"},
    @{Cell="D41"; Old="0.006899"; New="0.006890"},
    @{Cell="D42"; Old="0.1183"; New="0.1182"},
    @{Cell="D43"; Old="0.003318"; New="0.003454"},
    @{Cell="D44"; Old="0.01248"; New="0.01247"},
    @{Cell="D45"; Old="0.00006266"; New="0.00006243"},
    @{Cell="D46"; Old="0.00000000748"; New="0.00000000749"},
    @{Cell="D47"; Old="0.6978"; New="0.6982"},
    @{Cell="D48"; Old="0.1232"; New="0.1245"},
    @{Cell="D49"; Old="0.00002093"; New="0.00002096"},
    @{Cell="D50"; Old="0.01236"; New="0.01238"}
)

foreach ($chg in $changes) {
    $cell = $ws.Range($chg.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $chg.New
    $cell.ClearFormats()
}
